# The sheet originally listed an example/placeholder row (row 2,
# "exemplo_nfse.pdf.pdf") followed by an empty "minha_nota.pdf" row (row 3)
# and the real processed invoice data (row 4, "nota_goiania.pdf").
#
# The update wires up the "choose layout and process" function: the
# placeholder file name on row 2 is replaced by the real source file name
# ("minha_nota.pdf"), and the now-redundant empty row 3 is removed so the
# processed invoice data shifts up to row 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the empty "minha_nota.pdf" row; this shifts the data row below
# (nota_goiania.pdf, row 4) up to row 3.
$ws.Rows.Item(3).Delete()

# Update the placeholder source file name on row 2 to the real file name.
$ws.Range("A2").Value = "minha_nota.pdf"
